$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($r = 3; $r -le 32; $r++) {
    $periodCell = $ws.Cells.Item($r, 8)   # Column H - PERIOD TO EXPIRE
    $periodCell.Value = $periodCell.Value() - 1

    $updateCell = $ws.Cells.Item($r, 9)   # Column I - LAST UPDATE
    # Leading apostrophe forces Excel to keep this as literal text instead
    # of auto-converting the date-looking string into a date serial value.
    $updateCell.Value = "'04-Nov-2025"
}
